# Normalize URNs before looking for duplicates in imports.
#
# The duplicate-urn.xlsx test fixture previously stored column A (URN:Tel)
# as plain numbers and used differently-cased spellings of the same name in
# column B to exercise duplicate detection. Update the fixture so that:
#   - column A holds the phone numbers as *text* (some with leading "+"
#     and/or punctuation/spacing) so the importer's URN-normalization path
#     gets exercised, and
#   - column B uses mixed/proper-case variants of the names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: +250788382382 / Eric Newcomer
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "+250788382382"
$ws.Range("B2").Value = "Eric Newcomer"

# Row 3: +250788383383 / Nic Pottier
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "+250788383383"
$ws.Range("B3").Value = "Nic Pottier"

# Row 4: (+250) 788 382382 / Jen Newcomer  (same number as row 2, written
# with punctuation/spacing, to test normalization-before-dedupe)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "(+250) 788 382382"
$ws.Range("B4").Value = "Jen Newcomer"
